{"js": "// Replace several pieces of text in the document body, matching the\n// unified diff:\n//   ZA-Information / Zentralarchiv f\u00fcr Empirische Sozialforschung -> ZUMA Nachrichten\n//   1997 -> 2009\n//   Meltdown -> Object Model Capabilities For Distributed Object Management.\n//   Moritz Lipp, ... Mike Hamburg -> Frank Manola\n//   RWTH Aachen -> Fachhochschule f\u00fcr Verwaltung und Dienstleistung (Altenholz, Reinfeld)\n//   Chausseestr. 29 -> Oldesloer Strasse 32\n\nconst replacements = [\n  [\n    \"ZA-Information / Zentralarchiv f\u00fcr Empirische Sozialforschung\",\n    \"ZUMA Nachrichten\",\n  ],\n  [\"1997\", \"2009\"],\n  [\"Meltdown\", \"Object Model Capabilities For Distributed Object Management.\"],\n  [\n    \"Moritz Lipp, Michael Schwarz , Daniel Gruss, Thomas Prescher , Werner Haas, Stefan Mangard, Paul Kocher, Daniel Genkin, Yuval Yarom, Mike Hamburg\",\n    \"Frank Manola\",\n  ],\n  [\n    \"RWTH Aachen\",\n    \"Fachhochschule f\u00fcr Verwaltung und Dienstleistung (Altenholz, Reinfeld)\",\n  ],\n  [\"Chausseestr. 29\", \"Oldesloer Strasse 32\"],\n];\n\nfor (const [findText, newText] of replacements) {\n  const results = context.document.body.search(findText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${findText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the text replacements described by the diff:\n#   ZA-Information / Zentralarchiv f\u00fcr Empirische Sozialforschung -> ZUMA Nachrichten\n#   1997 -> 2009\n#   Meltdown -> Object Model Capabilities For Distributed Object Management.\n#   Moritz Lipp, ... Mike Hamburg -> Frank Manola\n#   RWTH Aachen -> Fachhochschule f\u00fcr Verwaltung und Dienstleistung (Altenholz, Reinfeld)\n#   Chausseestr. 29 -> Oldesloer Strasse 32\n#\n# Each target string is the sole content of its own paragraph, so we match\n# paragraphs by their exact (trimmed) text and overwrite just that paragraph's\n# range -- this rewrites the <w:t> in place without disturbing the paragraph\n# mark or any other run/paragraph formatting.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"ZA-Information / Zentralarchiv f\u00fcr Empirische Sozialforschung\"; Replace = \"ZUMA Nachrichten\" },\n    @{ Find = \"1997\"; Replace = \"2009\" },\n    @{ Find = \"Meltdown\"; Replace = \"Object Model Capabilities For Distributed Object Management.\" },\n    @{ Find = \"Moritz Lipp, Michael Schwarz , Daniel Gruss, Thomas Prescher , Werner Haas, Stefan Mangard, Paul Kocher, Daniel Genkin, Yuval Yarom, Mike Hamburg\"; Replace = \"Frank Manola\" },\n    @{ Find = \"RWTH Aachen\"; Replace = \"Fachhochschule f\u00fcr Verwaltung und Dienstleistung (Altenholz, Reinfeld)\" },\n    @{ Find = \"Chausseestr. 29\"; Replace = \"Oldesloer Strasse 32\" }\n)\n\nforeach ($r in $replacements) {\n    $done = $false\n    foreach ($p in $d.Paragraphs) {\n        $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($t -eq $r.Find) {\n            $p.Range.Text = $r.Replace\n            $done = $true\n            break\n        }\n    }\n    if (-not $done) {\n        # Fallback: plain Find & Replace across the whole document in case a\n        # target string ever spans/mixes with other content in a paragraph.\n        $rng = $d.Content\n        $find = $rng.Find\n        $find.ClearFormatting()\n        $find.Replacement.ClearFormatting()\n        $find.Execute(\n            $r.Find,\n            $true,\n            $false,\n            $false,\n            $false,\n            $false,\n            $true,\n            1,\n            $false,\n            $r.Replace,\n            2\n        )\n    }\n}\n"}
